$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append new row 84, mirroring the formatting of the last existing data row (83):
#   - column A gets the date-formatted style used by all prior date cells
#   - column B gets the plain numeric value
$ws.Cells.Item(83, 1).Copy($ws.Cells.Item(84, 1))
$ws.Cells.Item(84, 1).Value = 45884

$ws.Cells.Item(84, 2).Value = 0.768168485846715
